$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# STEP 1: Rename the chapter-title bookmark (_Toc455431969 -> _Toc455582657)
#         Bookmarks whose name starts with "_Toc" are excluded from the
#         Document.Bookmarks collection (same as real Word), so rebuild
#         the paragraph's XML with the new bookmark name.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(2)
$bmXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="CHAPTERTITLE"/></w:pPr><w:bookmarkStart w:id="0" w:name="_Toc455582657"/><w:r><w:rPr><w:sz w:val="36"/><w:szCs w:val="36"/></w:rPr><w:t>CONCLUSION</w:t></w:r><w:bookmarkEnd w:id="0"/></w:p>'
$p2.Range.InsertXML($bmXml)

# ---------------------------------------------------------------------
# STEP 2: Paragraph "ALI was designed..." - superscript "o" -> "◦" (FOV of 6°)
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(4)
$rng = $p4.Range.Duplicate
$rng.Find.Execute("large FOV of 6", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$oRng = $d.Range($rng.End, $rng.End + 1)
$oRng.Text = "◦"

# ---------------------------------------------------------------------
# STEP 3: Reorganize paragraphs. The "With the completion of ALI..."
#         paragraph moves from position 5 down to the very end of this
#         list (after "For an additional improvement..."); everything
#         else shifts up by one.
# ---------------------------------------------------------------------
$pMove = $d.Paragraphs.Item(5)
$pMove.Range.Cut() | Out-Null

$pLast = $d.Paragraphs.Item(9)   # now "For an additional improvement..." after the cut
$pasteTarget = $d.Range($pLast.Range.End, $pLast.Range.End)
$pasteTarget.Paste()

# ---------------------------------------------------------------------
# STEP 4: Degree-symbol fixes inside the paragraph that is now #5
#         ("The test flight for ALI occurred...")
# ---------------------------------------------------------------------
$pFlight = $d.Paragraphs.Item(5)

# "...SSA being relatively close to 90<o>" -> "...90<◦>"
$r = $pFlight.Range.Duplicate
$r.Find.Execute("relatively close to 90", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$supRng = $d.Range($r.End, $r.End + 1)
$supRng.Text = "◦"

# "... at 98<o >which" -> "... at 98<◦ >which"
$r2 = $pFlight.Range.Duplicate
$r2.Find.Execute(" at 98", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$supRng2 = $d.Range($r2.End, $r2.End + 2)
$supRng2.Text = "◦ "

# ---------------------------------------------------------------------
# STEP 5: Degree-symbol fixes inside the paragraph that is now #6
#         ("This first prototype ALI instrument...")
# ---------------------------------------------------------------------
$pProto = $d.Paragraphs.Item(6)

# "...SAA is in between 45-60<o>" -> "45-60<◦>"
$r3 = $pProto.Range.Duplicate
$r3.Find.Execute("in between 45-60", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$supRng3 = $d.Range($r3.End, $r3.End + 1)
$supRng3.Text = "◦"

# "...scattering angles of 80-100<o>" -> "80-100<◦>"
$r4 = $pProto.Range.Duplicate
$r4.Find.Execute("scattering angles of 80-100", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$supRng4 = $d.Range($r4.End, $r4.End + 1)
$supRng4.Text = "◦"

# ---------------------------------------------------------------------
# STEP 6: Text tweak in the paragraph that is now #10
#         ("With the completion of ALI..."), which was moved in step 3.
#         Add "for a space mission" and drop the trailing space.
# ---------------------------------------------------------------------
$pSim = $d.Paragraphs.Item(10)

$r5 = $pSim.Range.Duplicate
$r5.Find.Execute("over the total radiance. Overall", $false, $false, $false, $false, $false, $true, 1, $false, "over the total radiance for a space mission. Overall", 2) | Out-Null

# Remove the trailing space just before the final paragraph mark.
$tailStart = $pSim.Range.End - 2
$tailRng = $d.Range($tailStart, $tailStart + 1)
if ($tailRng.Text -eq " ") {
    $tailRng.Text = ""
}

# ---------------------------------------------------------------------
# STEP 7: Footer / header page numbers.
#         This section only defines a "first page" footer (wdHeaderFooterFirstPage,
#         index 2) and a "default"/primary header (wdHeaderFooterPrimary, index 1).
# ---------------------------------------------------------------------
$sec1 = $d.Sections.Item(1)
$footer1 = $sec1.Footers.Item(2)
$footer1.Range.Find.Execute("165", $false, $false, $false, $false, $false, $true, 1, $false, "168", 2) | Out-Null

$header1 = $sec1.Headers.Item(1)
$header1.Range.Find.Execute("167", $false, $false, $false, $false, $false, $true, 1, $false, "171", 2) | Out-Null
